$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5836.9375
$ws.Range("I74").Value = 5425.25
$ws.Range("K74").Value = 5425.25
$ws.Range("M74").Value = -4489.25

$ws.Range("H77").Value = 5836.9375
$ws.Range("I77").Value = 5425.25
$ws.Range("K77").Value = 27126.25
$ws.Range("M77").Value = -22446.25

$ws.Range("H114").Value = 40200
$ws.Range("J114").Value = 40200
$ws.Range("L114").Value = 40200
$ws.Range("N114").Value = -48878

$ws.Range("H127").Value = 2377.8572
$ws.Range("J127").Value = 3475
$ws.Range("L127").Value = 10425
$ws.Range("N127").Value = -20345

$ws.Range("H133").Value = 139000
$ws.Range("J133").Value = 139000
$ws.Range("L133").Value = 139000
$ws.Range("N133").Value = -149120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4451.544
$ws.Range("I32").Value = 4154.407
$ws.Range("K32").Value = 4154.407
$ws.Range("M32").Value = -3867.407

$ws.Range("H74").Value = 35558.027
$ws.Range("I74").Value = 45263.85
$ws.Range("J74").Value = 2800.875
$ws.Range("K74").Value = 45263.85
$ws.Range("L74").Value = 2800.875
$ws.Range("M74").Value = -44389.85
$ws.Range("N74").Value = -4548.875

$ws.Range("H76").Value = 80000
$ws.Range("J76").Value = 80000
$ws.Range("L76").Value = 80000
$ws.Range("N76").Value = -80676

$ws.Range("H77").Value = 35558.027
$ws.Range("I77").Value = 45263.85
$ws.Range("J77").Value = 2800.875
$ws.Range("K77").Value = 226319.25
$ws.Range("L77").Value = 14004.375
$ws.Range("M77").Value = -221951.25
$ws.Range("N77").Value = -22740.375

$ws.Range("H79").Value = 80000
$ws.Range("J79").Value = 80000
$ws.Range("L79").Value = 80000
$ws.Range("N79").Value = -82340

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 3059.4285
$ws.Range("I37").Value = 3157.4
$ws.Range("J37").Value = 2814.5
$ws.Range("K37").Value = 3157.4
$ws.Range("L37").Value = 2814.5
$ws.Range("M37").Value = -3020.4
$ws.Range("N37").Value = -3088.5

$ws.Range("H54").Value = 3366.3333
$ws.Range("J54").Value = 4000
$ws.Range("L54").Value = 4000
$ws.Range("N54").Value = -4968

$ws.Range("H134").Value = 1221
$ws.Range("I134").Value = 1221
$ws.Range("K134").Value = 3663
$ws.Range("M134").Value = -1128

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3692.1082
$ws.Range("I132").Value = 3459.2424
$ws.Range("K132").Value = 10377.7272
$ws.Range("M132").Value = -7847.727200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 37887.25
$ws.Range("J9").Value = 400
$ws.Range("L9").Value = 1200
$ws.Range("N9").Value = -1648

$ws.Range("H32").Value = 93306.82000000001
$ws.Range("I32").Value = 113375
$ws.Range("K32").Value = 340125
$ws.Range("M32").Value = -339842

$ws.Range("H110").Value = 5333.3335
$ws.Range("I110").Value = 5333.3335
$ws.Range("K110").Value = 16000.0005
$ws.Range("M110").Value = -11910.0005

$ws.Range("H131").Value = 995.2258
$ws.Range("I131").Value = 741.1429000000001
$ws.Range("J131").Value = 1204.4706
$ws.Range("K131").Value = 2223.4287
$ws.Range("L131").Value = 3613.4118
$ws.Range("M131").Value = 2816.5713
$ws.Range("N131").Value = -13693.4118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 22000
$ws.Range("I18").Value = 22000
$ws.Range("K18").Value = 22000
$ws.Range("M18").Value = -21707

$ws.Range("H21").Value = 15495.875
$ws.Range("I21").Value = 13994.667
$ws.Range("J21").Value = 19999.5
$ws.Range("K21").Value = 13994.667
$ws.Range("L21").Value = 19999.5
$ws.Range("M21").Value = -13821.667
$ws.Range("N21").Value = -20345.5

$ws.Range("H30").Value = 15495.875
$ws.Range("I30").Value = 13994.667
$ws.Range("J30").Value = 19999.5
$ws.Range("K30").Value = 13994.667
$ws.Range("L30").Value = 19999.5
$ws.Range("M30").Value = -13889.667
$ws.Range("N30").Value = -20209.5

$ws.Range("H35").Value = 22499.5
$ws.Range("I35").Value = 22499.5
$ws.Range("K35").Value = 22499.5
$ws.Range("M35").Value = -22201.5

$ws.Range("H122").Value = 3458.2
$ws.Range("J122").Value = 4007.5
$ws.Range("L122").Value = 12022.5
$ws.Range("N122").Value = -16922.5

$ws.Range("H132").Value = 779.3333
$ws.Range("J132").Value = 978.5
$ws.Range("L132").Value = 2935.5
$ws.Range("N132").Value = -7995.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1127.8462
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2340

$ws.Range("H40").Value = 2505.2222
$ws.Range("I40").Value = 2224
$ws.Range("J40").Value = 2856.75
$ws.Range("K40").Value = 2224
$ws.Range("L40").Value = 2856.75
$ws.Range("M40").Value = -2088
$ws.Range("N40").Value = -3128.75

$ws.Range("H43").Value = 22661.111

$ws.Range("H50").Value = 41746.25

$ws.Range("H69").Value = 48333
$ws.Range("I69").Value = 44999
$ws.Range("J69").Value = 50000
$ws.Range("K69").Value = 44999
$ws.Range("L69").Value = 50000
$ws.Range("M69").Value = -44188
$ws.Range("N69").Value = -51622

$ws.Range("H72").Value = 48333
$ws.Range("I72").Value = 44999
$ws.Range("J72").Value = 50000
$ws.Range("K72").Value = 134997
$ws.Range("L72").Value = 150000
$ws.Range("M72").Value = -130941
$ws.Range("N72").Value = -158112

$ws.Range("H132").Value = 2501.5
$ws.Range("I132").Value = 2290.3333
$ws.Range("J132").Value = 3768.5
$ws.Range("K132").Value = 6870.999899999999
$ws.Range("L132").Value = 11305.5
$ws.Range("M132").Value = -4340.999899999999
$ws.Range("N132").Value = -16365.5

$ws.Range("H133").Value = 68456
$ws.Range("J133").Value = 68456
$ws.Range("L133").Value = 68456
$ws.Range("N133").Value = -73516

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 25000
$ws.Range("I43").Value = 24791.666
$ws.Range("K43").Value = 24791.666
$ws.Range("M43").Value = -24642.666

$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws.Range("H136").Value = 3172.4355
$ws.Range("I136").Value = 2119.1135
$ws.Range("J136").Value = 5747.222
$ws.Range("K136").Value = 6357.3405
$ws.Range("L136").Value = 17241.666
$ws.Range("M136").Value = -3807.3405
$ws.Range("N136").Value = -22341.666
